$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "44.443.92"
$ws.Range("E2").Value = "  +2.47%  "

# Row 3
$ws.Range("D3").Value = "2.366.56"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("E5").Value = "  +3.74%  "

# Row 6
Set-TextValue $ws.Range("D6") "239.09"
$ws.Range("E6").Value = "  +2.89%  "

# Row 7
Set-TextValue $ws.Range("D7") "73.48"
$ws.Range("E7").Value = "  +7.76%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.552"
$ws.Range("E9").Value = "  +20.08%  "

# Row 10
$ws.Range("E10").Value = "  +7.31%  "

# Row 11
Set-TextValue $ws.Range("D11") "29.95"
$ws.Range("E11").Value = "  +12.91%  "

# Row 12
$ws.Range("E12").Value = "  +2.23%  "

# Row 13
$ws.Range("D13").Value = "2.716.07"
$ws.Range("E13").Value = "  -0.07%  "

# Row 14
Set-TextValue $ws.Range("D14") "16.90"
$ws.Range("E14").Value = "  +8.03%  "

# Row 15
$ws.Range("E15").Value = "  +9.08%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.904"
$ws.Range("E16").Value = "  +7.43%  "

# Row 17
$ws.Range("D17").Value = "2.369.33"
$ws.Range("E17").Value = "  +0.15%  "

# Row 18
$ws.Range("D18").Value = "44.476.22"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19
$ws.Range("E19").Value = "  +5.02%  "

# Row 20
Set-TextValue $ws.Range("D20") "77.49"
$ws.Range("E20").Value = "  +4.70%  "

# Row 21
Set-TextValue $ws.Range("D21") "6.48"
$ws.Range("E21").Value = "  +4.06%  "

# Row 22
Set-TextValue $ws.Range("D22") "254.80"
$ws.Range("E22").Value = "  +2.69%  "

# Row 23
Set-TextValue $ws.Range("D23") "3.85"
$ws.Range("E23").Value = "  -2.45%  "

# Row 24
$ws.Range("E24").Value = "  -0.01%  "

# Row 25
$ws.Range("E25").Value = "  +2.54%  "

# Row 26
Set-TextValue $ws.Range("D26") "10.42"
$ws.Range("E26").Value = "  +4.49%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.30"
$ws.Range("E27").Value = "  +3.48%  "

# Row 28
Set-TextValue $ws.Range("D28") "22.50"
$ws.Range("E28").Value = "  +0.86%  "

# Row 29
$ws.Range("E29").Value = "  +4.49%  "

# Row 30
Set-TextValue $ws.Range("D30") "173.96"
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.131"
$ws.Range("E31").Value = "  +2.51%  "

# Row 32
$ws.Range("E32").Value = "  +5.41%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0744"
$ws.Range("E33").Value = "  +7.68%  "

# Row 34
$ws.Range("E34").Value = "  +4.50%  "

# Row 35
$ws.Range("E35").Value = "  +3.27%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.92"
$ws.Range("E36").Value = "  +7.76%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.43"
$ws.Range("E37").Value = "  -2.93%  "

# Row 38
Set-TextValue $ws.Range("D38") "6.49"
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("E39").Value = "  +6.85%  "

# Row 40
Set-TextValue $ws.Range("D40") "20.03"
$ws.Range("E40").Value = "  +10.26%  "

# Row 41
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.85"
$ws.Range("E42").Value = "  -1.14%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.25"
$ws.Range("E43").Value = "  +3.56%  "

# Row 44
$ws.Range("E44").Value = "  +3.93%  "

# Row 45
$ws.Range("E45").Value = "  +1.38%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D46") "0.185"
$ws.Range("E46").Value = "  +12.86%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "98.87"
$ws.Range("E47").Value = "  +0.08%  "

# Row 48
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D48") "4.48"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.35"
$ws.Range("E49").Value = "  +3.43%  "

# Row 50
$ws.Range("D50").Value = "1.445.18"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").Value = "2.590.43"
$ws.Range("E51").Value = "  -0.07%  "
